# Update the "Next Appointment" date and "Status" for the first few pending
# reminder rows, as part of the send_reminder workflow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: appointment date -> 2025-10-24, status -> Reminded
$ws.Range("C2").Value = 45954
$ws.Range("D2").Value = "Reminded"

# Row 3: appointment date -> 2025-10-24 (status already Reminded)
$ws.Range("C3").Value = 45954

# Row 4: appointment date -> 2025-10-24, status -> Reminded
$ws.Range("C4").Value = 45954
$ws.Range("D4").Value = "Reminded"

# Row 5: appointment date -> 2025-10-24, status -> Reminded
$ws.Range("C5").Value = 45954
$ws.Range("D5").Value = "Reminded"
